$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily rows to append, continuing the existing time series
# (date serial, nuovi pos., somma mobile 7gg., somma mobile 7gg. per 100mila abitanti)
$data = @(
    @(44330, 1, 6, 39.5882818685669),
    @(44331, 0, 6, 39.5882818685669),
    @(44332, 2, 8, 52.78437582475588),
    @(44333, 0, 7, 46.18632884666139),
    @(44334, 0, 6, 39.5882818685669),
    @(44335, 0, 6, 39.5882818685669),
    @(44336, 0, 3, 19.79414093428345),
    @(44337, 3, 5, 32.99023489047242),
    @(44338, 0, 5, 32.99023489047242),
    @(44339, 0, 3, 19.79414093428345),
    @(44340, 0, 3, 19.79414093428345),
    @(44341, 0, 3, 19.79414093428345),
    @(44342, 0, 3, 19.79414093428345),
    @(44343, 0, 3, 19.79414093428345)
)

$lastRow = 255
$startRow = $lastRow + 1

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]

    # copy the formatting (styles) of the last existing data row down to the new row
    $ws.Range("A" + $lastRow + ":D" + $lastRow).Copy()
    $ws.Range("A" + $r + ":D" + $r).PasteSpecial(-4122)

    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
